$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A13").Value = "W.A.901150038@mailsac.com"
$ws.Range("B13").Value = "Abcd1234"
